# Contribution sheet updated for week 4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Group 2"
$ws.Name = "Group 2"

# Fill in the Week 4 / "Code Review 4" column (F) contribution scores
$ws.Range("F2").Value = 25
$ws.Range("F3").Value = 25
$ws.Range("F4").Value = 25
$ws.Range("F5").Value = 25

# Move the active cell selection to F6, matching where the user left off
# after entering this week's figures
$ws.Range("F6").Select()
